$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (rows 2-67) from 45192 to 45202
$ws.Range("C2:C67").Value = 45202
